$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shortened NACE industry labels (column B), same 21-label cycle repeats for
# the Haugalandet block (rows 2-22) and the Norge block (rows 23-43).
$labels = @(
    "A - (01-03)",
    "B - (05-09)",
    "C - (10-33)",
    "D - (35)",
    "E - (36-39)",
    "F - (41-43)",
    "G - (45-47)",
    "H - (49-53)",
    "I - (55-56)",
    "J - (58-63)",
    "K - (64-66)",
    "L - (68)",
    "M - (69-75)",
    "N - (77-82)",
    "O - (84)",
    "P - (85)",
    "Q - (86-88)",
    "R - (90-93)",
    "S - (94-96)",
    "T - (97)",
    "U - (99)"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $labels[$i]
    $ws.Cells.Item(23 + $i, 2).Value = $labels[$i]
}

# Widen column B to fit the (still fairly long) labels.
$ws.Columns("B").ColumnWidth = 44.166666666666664

# Move the active selection to B47 (below the data, as left by the author).
$ws.Range("B47").Select()
